# Pandas_Cheat_Sheet.pptx - "Frequently Used Options" textbox update.
# - Italicised the `option_name` placeholder wherever it appears.
# - Cleaned up the `pd.options.option_name` / `pd.options.display.max_rows`
#   examples (removed stray underscores, fixed the stray capital P).
# - Reworded the last line describing the max_rows example.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item("TextBox 37")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# ---------------------------------------------------------------------
# Paragraph 2: "Options can be queried and set via: pd.options._option_name_
#               (where _option_name_ is the name of an option). For example:"
# ---------------------------------------------------------------------

# a) "pd.options._" -> "pd.options."
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.IndexOf("pd.options._")
$abs = $pStart + $rel
$tr.Characters($abs, ("pd.options._").Length).Text = "pd.options."

# b) remove the stray "_" that used to directly follow the first "option_name"
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.IndexOf("option_name_")
$abs = $pStart + $rel + ("option_name").Length
$tr.Characters($abs, 1).Text = ""

# c) italicise the first "option_name" (bold, Consolas run)
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.IndexOf("option_name")
$abs = $pStart + $rel
$tr.Characters($abs, ("option_name").Length).Font.Italic = -1

# d) " (where _" -> " (where "
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.IndexOf(" (where _")
$abs = $pStart + $rel
$tr.Characters($abs, (" (where _").Length).Text = " (where "

# e) italicise the second "option_name" (plain run)
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.LastIndexOf("option_name")
$abs = $pStart + $rel
$tr.Characters($abs, ("option_name").Length).Font.Italic = -1

# f) "_ is the name of an option). For example:" -> " is the name of an option). For example:"
$para2 = $tr.Paragraphs(2, 1)
$pStart = $para2.Start
$pText = $para2.Text
$rel = $pText.IndexOf("_ is the name of an option)")
$abs = $pStart + $rel
$tr.Characters($abs, 1).Text = ""

# ---------------------------------------------------------------------
# Paragraph 3: "Pd.options.display.max_rows = 20"
# ---------------------------------------------------------------------
$para3 = $tr.Paragraphs(3, 1)
$pStart = $para3.Start
$pText = $para3.Text
$rel = $pText.IndexOf("Pd.options.display.max_rows")
$abs = $pStart + $rel
$tr.Characters($abs, ("Pd.options.display.max_rows").Length).Text = "pd.options.display.max_rows"

# ---------------------------------------------------------------------
# Paragraph 4: "'max_rows' option is currently set to 20"
#           -> "Set the display.max_rows option to 20."
# ---------------------------------------------------------------------
$para4 = $tr.Paragraphs(4, 1)
$para4.Text = "Set the display.max_rows option to 20."

$para4 = $tr.Paragraphs(4, 1)
$pStart = $para4.Start
$pText = $para4.Text
$rel = $pText.IndexOf("display.max_rows")
$abs = $pStart + $rel
$fmtRun = $tr.Characters($abs, ("display.max_rows").Length)
$fmtRun.Font.Bold = -1
$fmtRun.Font.Name = "Consolas"
